$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix borders on rows 45-46 (columns A-C) -----------------------------
# Before this edit, A45:C46 used the border-less style (s="3"). The edit adds
# the thin border used everywhere else in the table (matching the style
# already used by D45, which has font3 + border1, no wrap/alignment override
# once WrapText is turned back off).
$ws.Range("D45").Copy()
$ws.Range("A45:C46").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A45:C46").WrapText = $false
$excel.CutCopyMode = $false

# --- 2. Append new row 47 with the new "2 Oklar ayrı paralel Dot" project --
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "DP-000-2D-NUV-LIP-P2B1-01"
$ws.Range("C47").Value = "Bağımsız Tasarım"
$ws.Range("D47").Value = "2 Oklu 3x Dot"
$ws.Range("E47").Value = "Nuvo Mcu"
$ws.Range("F47").Value = "Kapı Üstü"
$ws.Range("G47").Value = "Paralel"
$ws.Range("H47").Value = "Kablolu"
$ws.Range("I47").Value = "Buzzerlı"
$ws.Range("J47").Value = "Model-01"
$ws.Range("K47").Value = "https://github.com/btk42/DP-000-2D-NUV-LIP-P2B1-01"

# Row height to match the rest of the table
$ws.Rows.Item(47).RowHeight = 21

# Hyperlink for the new row's Link cell (do this before the format paste
# below so the final cell format matches the other Link-column cells)
$ws.Hyperlinks.Add($ws.Range("K47"), "https://github.com/btk42/DP-000-2D-NUV-LIP-P2B1-01") | Out-Null

# Copy matching formats down from the row above
$ws.Range("A45").Copy()
$ws.Range("A47").PasteSpecial(-4122)
$ws.Range("B45").Copy()
$ws.Range("B47").PasteSpecial(-4122)
$ws.Range("D45").Copy()
$ws.Range("C47:J47").PasteSpecial(-4122)
$ws.Range("K45").Copy()
$ws.Range("K47").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Selection / view ------------------------------------------------
$ws.Range("M36").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
